# This commit replaces the PPTI class-boundary lookup table on Sheet1.
# The old 14 "LN*"/"L30x" type rows are replaced by a new, larger table of
# 25 "L1xx"/"L2xx"/"L3xx" type rows (rows 2-26); the L301-L306 rows are kept
# (with identical values) at the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rowsData = @(
        @("L107",5,4,2.86,2.6,2.4300000000000002,2.2599999999999998,2.09,1.1599999999999999),
        @("L109",5,4,2.86,2.6,2.4300000000000002,2.2599999999999998,2.09,1.1599999999999999),
        @("L104",5,4,2.69,2.5099999999999998,2.34,2.17,2,1.1599999999999999),
        @("L105a",5,4,2.69,2.5099999999999998,2.34,2.17,2,1.1599999999999999),
        @("L207",5,4,2.69,2.5099999999999998,2.34,2.17,2,1.1599999999999999),
        @("L105b",5,4,2.6,2.4300000000000002,2.2599999999999998,2.09,1.9,1.1599999999999999),
        @("L106",5,4,2.86,2.6,2.4300000000000002,2.2599999999999998,2.09,1.1599999999999999),
        @("L208",5,4,2.86,2.6,2.4300000000000002,2.2599999999999998,2.09,1.1599999999999999),
        @("L101",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L102",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L201",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L202",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L204",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L205",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L103",5,4,2.69,2.5099999999999998,2.34,2.17,2,1.1599999999999999),
        @("L203",5,4,2.69,2.5099999999999998,2.34,2.17,2,1.1599999999999999),
        @("L206",5,4,2.69,2.5099999999999998,2.34,2.17,2,1.1599999999999999),
        @("L108",5,4,3.07,2.73,2.56,2.39,2.2200000000000002,1.1599999999999999),
        @("L110",5,4,3.07,2.73,2.56,2.39,2.2200000000000002,1.1599999999999999),
        @("L301",5,4,2.41,2.2400000000000002,2.0699999999999998,1.9,1.7,1.1599999999999999),
        @("L302",5,4,2.41,2.2400000000000002,2.0699999999999998,1.9,1.7,1.1599999999999999),
        @("L304",5,4,2.41,2.2400000000000002,2.0699999999999998,1.9,1.7,1.1599999999999999),
        @("L305",5,4,2.41,2.2400000000000002,2.0699999999999998,1.9,1.7,1.1599999999999999),
        @("L303",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999),
        @("L306",5,4,2.5099999999999998,2.34,2.17,2,1.8,1.1599999999999999)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowArr = $rowsData[$i]
    for ($j = 0; $j -lt $rowArr.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowArr[$j]
    }
}

$wb.Save()
